# DPLKAKT062-001 - Setup Periode Harian - Tambah Data.xlsx
# Commit: "Update Regresi Tanggal 31/03/2023"
#
# The underlying edit: the TGL_TRANS ("Tgl Trans") value in cell O2 is
# updated from "05/05/2023" to "07/09/2024". The dependent PERIODE_HARIAN
# formula in P2 (=RIGHT(O2,4) & MID(O2,4,2) & LEFT(O2,2)) recalculates
# automatically. The active selection on the sheet also moves to R2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the transaction date cell.
$ws.Range("O2").Value = "07/09/2024"

# Move the live selection to match the saved view (R2).
$ws.Range("R2").Select()
